$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("aragon") dimension metadata is re-curated into a reference-area
# dimension: the iaest-specific dimension URI is replaced by the generic
# sdmx refArea dimension, and the concept type becomes a URI to the new
# "Comunidad" (autonomous community) code list instead of a skos:Concept.
# The per-column mapping workbook entry (row 5) no longer applies and is
# cleared.
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E4").Value = "URI-Comunidad"
$ws.Range("E5").Clear()

# Column G ("combustible") moves from being a curated dimension to being a
# measure: iaest-dimension -> iaest-measure, "dim" -> "medida", and the
# concept type becomes a plain integer (xsd:int) rather than a
# skos:Concept. Its mapping workbook entry (row 5) is likewise cleared.
$ws.Range("G2").Value = "iaest-measure:combustible"
$ws.Range("G3").Value = "medida"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("G5").Clear()
